$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The refreshed data set drops two samples entirely: the original row 26
# ("RM 232") and the original row 28 ("SC 92"). Deleting row 26 first shifts
# "SC 92" up into row 27, so it is the next one removed - every row below
# slides up two spots and the used range shrinks from F35 to F33.
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# --- Column C ("B") re-imputation among the unchanged RM rows (2-25) ---
$ws.Range("C2").Value = 14.9
$ws.Range("C6").Value = ""
$ws.Range("C12").Value = 12.5
$ws.Range("C14").Value = ""
$ws.Range("C20").Value = 12.5
$ws.Range("C21").Value = 12.7
$ws.Range("C23").Value = ""
$ws.Range("C24").Value = ""

# --- Rows 26-33 (the SC rows, after the deletions above) get refreshed
#     values too, for columns B ("A") and/or C ("B") ---
$ws.Range("B26").Value = -20.2   # SC 5
$ws.Range("B27").Value = ""      # SC 101
$ws.Range("B30").Value = -19.7   # SC 120
$ws.Range("C31").Value = 15.3    # SC 132
$ws.Range("B32").Value = ""      # SC 193
$ws.Range("C33").Value = 10.4    # SC 232
